$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update party name headers (row 1, columns B:U) to full descriptive names
$ws.Range("B1").Value = 'BANAAN - Better Seeking Alternatives than Doing Nothing in Apathy (Beter Alternatieven Nastreven Als Apathisch Nietsdoen/ - , BANAAN), known until  as Radical Reformers fighting for an upright Society (ROSSEM, Radikale Omvormers Strijders en Strubbelaars voor een Eerlijke Maatschappij)'
$ws.Range("C1").Value = 'CD&V - Christian Democratic and Flemish (Christen-Democratisch en Vlaams, CD&V), known until 29 September 2001 as Christian People’s Party (CVP, Christelijke Volkspartij)'
$ws.Range("D1").Value = 'CDH - Democrat Humanist Centre (Centre Démocrate Humaniste, CDH), known until 18 May 2002 as Christian Social Party (PSC, / Parti social chrétien)'
$ws.Range("E1").Value = 'ECOLO - EcoIogists (Ecologistes pour I’organisation de luttes originales, ECOLO)'
$ws.Range("F1").Value = 'FDF-PPW - Francophone Democratic Front (Front democratique francophone_Pan pour la Wallonie, FDF-PPW)'
$ws.Range("G1").Value = 'FN-NF - National Front (Nationaal Front/ Front National, FN-NF)'
$ws.Range("H1").Value = 'G - Green (Groen, G), known until 15 November 2003 as Live differently (AGALEV, Anders gaan leven), known from 15 November 2003 until 12 January 2012 as Green! (G!, Groen!)'
$ws.Range("I1").Value = 'PRL - Liberal Reform Party (Parti réformateur liberal, PRL)'
$ws.Range("J1").Value = 'PS - Socialist Party (Parti socialiste, PS)'
$ws.Range("K1").Value = 'PVV - Party of Liberty and Progress ( Partij voor Vrijheid en Vooruitgang, PVV)'
$ws.Range("L1").Value = 'VB - Flemish Block (Vlaams Blok, VB)'
$ws.Range("M1").Value = 'VU - People''s Union (Volksunie, VU)'
$ws.Range("N1").Value = 'sp.a - Socialist Party-Differently (Socialistische Partij Anders, sp.a), known until 13 October 2001 as Socialist Party (SP, SocialistischePartij)'
$ws.Range("O1").Value = 'Open VLD - Open Flemish Liberals and Democrats (Open VLD) (Open Vlaamse Liberalen Demokraten , Open VLD)'
$ws.Range("P1").Value = 'MR - Reform Movement ( Mouvement Réformateur, MR)'
$ws.Range("Q1").Value = 'LDD - De Decker''s List (Lijst De Decker, LDD)'
$ws.Range("R1").Value = 'N-VA - New Flemish Alliance (Nieuw-Vlaams Alliantie, N-VA)'
$ws.Range("S1").Value = 'PP - People''s Party (Parti Populaire, PP)'
$ws.Range("T1").Value = 'PTB - Labour Party (Parti du Travail de Belgique, PTB)'
$ws.Range("U1").Value = 'PVDA-PTB - Labour Party (Partij van de Arbeid/Parti du Travail de Belgique, PVDA-PTB)'

# Clean up floating point rounding noise in the numeric seat-count cells
$ws.Range("B2").Value = 3
$ws.Range("D2").Value = 18
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 3
$ws.Range("H2").Value = 7
$ws.Range("N2").Value = 28
$ws.Range("C3").Value = 29
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 6
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 18
$ws.Range("J3").Value = 21
$ws.Range("L3").Value = 11
$ws.Range("M3").Value = 5
$ws.Range("N3").Value = 20
$ws.Range("O3").Value = 21
$ws.Range("C4").Value = 22
$ws.Range("E4").Value = 11
$ws.Range("H4").Value = 9
$ws.Range("J4").Value = 19
$ws.Range("M4").Value = 8
$ws.Range("N4").Value = 14
$ws.Range("O4").Value = 23
$ws.Range("C6").Value = 30
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 8
$ws.Range("H6").Value = 4
$ws.Range("J6").Value = 20
$ws.Range("L6").Value = 17
$ws.Range("N6").Value = 14
$ws.Range("O6").Value = 18
$ws.Range("P6").Value = 23
$ws.Range("Q6").Value = 5
$ws.Range("D7").Value = 9
$ws.Range("E7").Value = 8
$ws.Range("H7").Value = 5
$ws.Range("J7").Value = 26
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 2
$ws.Range("H9").Value = 8
$ws.Range("J9").Value = 20
$ws.Range("L9").Value = 18
$ws.Range("N9").Value = 9
$ws.Range("P9").Value = 14
$ws.Range("U9").Value = 12

Write-Output "edit complete"
